$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20: "Seguimiento de cadena" -> "Seguimiento de producciones"
$ws.Range("C20").Value = "Seguimiento de producciones"

# Row 25: clear the "?" / long explanatory comment, mark as "Pasa" instead
$ws.Range("H25").Value = "Pasa"
$ws.Range("I25").ClearContents()

# Row 26: fill in the redirection check that used to live on the (now removed) row 27
$ws.Range("G26").Value = "Redirección correcta a fabricaComenzarProduccion.html"

# Remove the old "Botón 11.1" / "Botón 11.2" rows (Crear producción Pilser/Stout) -
# their content has been consolidated into row 26 above.
$ws.Rows("27:28").Delete()

# Restore the view: scroll position + active selection
$window = $excel.ActiveWindow
$window.ScrollRow = 18
$window.ScrollColumn = 1
$ws.Range("G29").Select()
